# Demo1.docx -- "comiited in US time zone"
#
# Appends, after the existing "Demo1" paragraph:
#   (blank line)
#   Adding to word file now
#   (blank line)
#   Mind changed
#
# The blank lines need to end up as genuinely empty paragraphs (no run at
# all), which is what Word itself produces for a bare blank line. A plain
# Selection.TypeParagraph()/Range.InsertParagraphAfter() at the very end of
# the story leaves behind a placeholder empty run, so instead we type a
# one-off marker character where each blank line belongs, split it into its
# own paragraph (with real text on both sides, so nothing needs to fabricate
# a placeholder run), and then delete just the marker -- leaving a clean,
# empty paragraph behind.

$d = $word.ActiveDocument

# Marker used purely as a splice point for blank lines; Start-of-Heading
# (0x01) can't occur in normal prose so Find will only ever match our own
# insertions.
$marker = [char]1

$sel = $word.Selection
$sel.EndKey(6) | Out-Null   # wdStory -> jump to the very end of the document
$sel.TypeText($marker + "Adding to word file now" + $marker + "Mind changed")

function Split-MarkerIntoBlankParagraph {
    param($doc, $markerChar)

    $findRange = $doc.Content
    $found = $findRange.Find.Execute($markerChar, $false, $false, $false, $false, $false, `
                                      $true, 1, $false, $null, 0)
    if (-not $found) { return }

    $start = $findRange.Start
    $end = $findRange.End

    # Give the marker its own paragraph: break after it, then break before it.
    $doc.Range($end, $end).InsertParagraphAfter()
    $doc.Range($start, $start).InsertParagraphAfter()

    # The "break before" shifted the (now isolated) marker forward by one.
    # Deleting it leaves a bare, run-less paragraph -- a true blank line.
    $doc.Range($start + 1, $end + 1).Delete()
}

Split-MarkerIntoBlankParagraph $d $marker
Split-MarkerIntoBlankParagraph $d $marker
